# Trade #17 (leadlag, UP) recorded at 2026-02-16 21:58:33 - currently OPEN, +0.000%.
# Append the new trade as the last row of both the "All Trades" sheet and the
# strategy-specific "leadlag" sheet.

$wb = $excel.ActiveWorkbook

$tradeNum    = 17
# Leading apostrophe forces the date to be stored as literal text (matching
# every other Date cell in these sheets) instead of being auto-parsed into a
# real Excel date serial number.
$tradeDate   = "'2026-02-16"
$tradeTime   = "21:58:33"
$strategy    = "leadlag"
$side        = "UP"
$entryPrice  = 68429.955
$exitPrice   = ""
$status      = "OPEN"
$pnlPct      = 0
$pnlDollar   = 0
$capAfter    = 100.0660986376279
$confidence  = 0.75
$entryReason = "Coinbase leading with 0.078% move"
$exitReason  = ""
$duration    = 0

function Add-TradeRow($ws, $row) {
    $ws.Cells.Item($row, 1).Value  = $tradeNum
    $ws.Cells.Item($row, 2).Value  = $tradeDate
    $ws.Cells.Item($row, 3).Value  = $tradeTime
    $ws.Cells.Item($row, 4).Value  = $strategy
    $ws.Cells.Item($row, 5).Value  = $side
    $ws.Cells.Item($row, 6).Value  = $entryPrice
    $ws.Cells.Item($row, 7).Value  = $exitPrice
    $ws.Cells.Item($row, 8).Value  = $status
    $ws.Cells.Item($row, 9).Value  = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capAfter
    $ws.Cells.Item($row, 12).Value = $confidence
    $ws.Cells.Item($row, 13).Value = $entryReason
    $ws.Cells.Item($row, 14).Value = $exitReason
    $ws.Cells.Item($row, 15).Value = $duration
}

# "All Trades" sheet: existing data occupies rows 1-17 (header + 16 trades),
# the new trade becomes row 18.
$wsAll = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAll 18

# "leadlag" sheet: existing data occupies rows 1-16 (header + 15 trades),
# the new trade becomes row 17.
$wsLeadLag = $wb.Worksheets.Item("leadlag")
Add-TradeRow $wsLeadLag 17
